# Fix bug: Carnaval event end date incorrectly differed from the start date.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("eventosSazonais")

# Correct the start (B2) and end (C2) dates for the "Carnaval" row so both
# point to the same date (13/02/2023 -> serial 44970).
$ws.Range("B2").Value = 44970
$ws.Range("C2").Value = 44970

# Move the active selection to C3, matching the saved cursor position.
$ws.Range("C3").Select()
